$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 814.0909
$ws.Range("I2").Value = 1001.875
$ws.Range("J2").Value = 313.33334
$ws.Range("K2").Value = 1001.875
$ws.Range("L2").Value = 313.33334
$ws.Range("M2").Value = -888.875
$ws.Range("N2").Value = -539.33334
$ws.Range("H13").Value = 48007.39
$ws.Range("J13").Value = 66.31579000000001
$ws.Range("L13").Value = 66.31579000000001
$ws.Range("N13").Value = -404.31579
$ws.Range("H43").Value = 4255.1113
$ws.Range("J43").Value = 5199.6665
$ws.Range("L43").Value = 5199.6665
$ws.Range("N43").Value = -5337.6665
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H132").Value = 1344.3954
$ws.Range("I132").Value = 1310.95
$ws.Range("K132").Value = 3932.85
$ws.Range("M132").Value = -1402.85
$ws.Range("H135").Value = 542.913
$ws.Range("I135").Value = 431.27274
$ws.Range("K135").Value = 3881.45466
$ws.Range("M135").Value = -1346.45466
$ws.Range("H137").Value = 2965.6177
$ws.Range("I137").Value = 2183.6736
$ws.Range("K137").Value = 6551.0208
$ws.Range("M137").Value = -4001.0208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3148.611
$ws.Range("I2").Value = 3104.2144
$ws.Range("J2").Value = 3304
$ws.Range("K2").Value = 3104.2144
$ws.Range("L2").Value = 3304
$ws.Range("M2").Value = -2991.2144
$ws.Range("N2").Value = -3530
$ws.Range("H32").Value = 22845.967
$ws.Range("I32").Value = 4698.1816
$ws.Range("K32").Value = 4698.1816
$ws.Range("M32").Value = -4411.1816
$ws.Range("H116").Value = 3148.611
$ws.Range("I116").Value = 3104.2144
$ws.Range("J116").Value = 3304
$ws.Range("K116").Value = 3104.2144
$ws.Range("L116").Value = 3304
$ws.Range("M116").Value = -810.2143999999998
$ws.Range("N116").Value = -7892
$ws.Range("H132").Value = 5827
$ws.Range("I132").Value = 3373.7
$ws.Range("J132").Value = 9915.833000000001
$ws.Range("K132").Value = 10121.1
$ws.Range("L132").Value = 29747.499
$ws.Range("M132").Value = -7591.099999999999
$ws.Range("N132").Value = -34807.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3148.611
$ws.Range("I3").Value = 3104.2144
$ws.Range("J3").Value = 3304
$ws.Range("K3").Value = 3104.2144
$ws.Range("L3").Value = 3304
$ws.Range("M3").Value = -2990.2144
$ws.Range("N3").Value = -3532
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H75").Value = 11165.667
$ws.Range("I75").Value = 11165.667
$ws.Range("K75").Value = 11165.667
$ws.Range("M75").Value = -10229.667
$ws.Range("H78").Value = 11165.667
$ws.Range("I78").Value = 11165.667
$ws.Range("K78").Value = 11165.667
$ws.Range("M78").Value = -28817.001
$ws.Range("H88").Value = 51234.3
$ws.Range("J88").Value = 51234.3
$ws.Range("L88").Value = 51234.3
$ws.Range("N88").Value = -52046.3
$ws.Range("H91").Value = 51234.3
$ws.Range("J91").Value = 51234.3
$ws.Range("L91").Value = 51234.3
$ws.Range("N91").Value = -54042.3
$ws.Range("H134").Value = 1595.965
$ws.Range("I134").Value = 1384.3585
$ws.Range("K134").Value = 4153.0755
$ws.Range("M134").Value = -1618.0755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1320141.2
$ws.Range("I6").Value = 1508590.1
$ws.Range("K6").Value = 1508590.1
$ws.Range("M6").Value = -1508477.1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H31").Value = 21312.057
$ws.Range("I31").Value = 32402.213
$ws.Range("K31").Value = 32402.213
$ws.Range("M31").Value = -32107.213
$ws.Range("H34").Value = 21312.057
$ws.Range("I34").Value = 32402.213
$ws.Range("K34").Value = 32402.213
$ws.Range("M34").Value = -32200.213
$ws.Range("H74").Value = 57624.2
$ws.Range("J74").Value = 66474.75
$ws.Range("L74").Value = 66474.75
$ws.Range("N74").Value = -68222.75
$ws.Range("H77").Value = 57624.2
$ws.Range("J77").Value = 66474.75
$ws.Range("L77").Value = 199424.25
$ws.Range("N77").Value = -208160.25
$ws.Range("H99").Value = 31229.285
$ws.Range("I99").Value = 55961
$ws.Range("J99").Value = 12680.5
$ws.Range("K99").Value = 55961
$ws.Range("L99").Value = 12680.5
$ws.Range("M99").Value = -54463
$ws.Range("N99").Value = -15676.5
$ws.Range("H122").Value = 40353.742
$ws.Range("I122").Value = 49864.57
$ws.Range("K122").Value = 149593.71
$ws.Range("M122").Value = -147143.71
$ws.Range("H126").Value = 31229.285
$ws.Range("I126").Value = 55961
$ws.Range("J126").Value = 12680.5
$ws.Range("K126").Value = 167883
$ws.Range("L126").Value = 38041.5
$ws.Range("M126").Value = -165413
$ws.Range("N126").Value = -42981.5
$ws.Range("H132").Value = 2293.509
$ws.Range("I132").Value = 2212.761
$ws.Range("J132").Value = 2706.2222
$ws.Range("K132").Value = 6638.282999999999
$ws.Range("L132").Value = 8118.6666
$ws.Range("M132").Value = -4108.282999999999
$ws.Range("N132").Value = -13178.6666
$ws.Range("H134").Value = 1920.6552
$ws.Range("I134").Value = 1790.2653
$ws.Range("J134").Value = 2630.5557
$ws.Range("K134").Value = 5370.7959
$ws.Range("L134").Value = 7891.6671
$ws.Range("M134").Value = -2835.7959
$ws.Range("N134").Value = -12961.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3334076.8
$ws.Range("I11").Value = 231
$ws.Range("J11").Value = 5000999.5
$ws.Range("K11").Value = 693
$ws.Range("L11").Value = 15002998.5
$ws.Range("M11").Value = -553
$ws.Range("N11").Value = -15003278.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H130").Value = 64333.332
$ws.Range("J130").Value = 64333.332
$ws.Range("L130").Value = 64333.332
$ws.Range("N130").Value = -74373.33199999999
$ws.Range("H132").Value = 3047.9744
$ws.Range("I132").Value = 2965.818
$ws.Range("J132").Value = 3499.8333
$ws.Range("K132").Value = 8897.454000000002
$ws.Range("L132").Value = 10499.4999
$ws.Range("M132").Value = -6367.454000000002
$ws.Range("N132").Value = -15559.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 19791.438
$ws.Range("I7").Value = 35717.875
$ws.Range("K7").Value = 35717.875
$ws.Range("M7").Value = -35605.875
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H46").Value = 25198.611
$ws.Range("I46").Value = 61617.855
$ws.Range("J46").Value = 2022.7273
$ws.Range("K46").Value = 61617.855
$ws.Range("L46").Value = 2022.7273
$ws.Range("M46").Value = -61429.855
$ws.Range("N46").Value = -2398.7273
$ws.Range("H61").Value = 75043.94500000001
$ws.Range("I61").Value = 67430.664
$ws.Range("J61").Value = 113110.336
$ws.Range("K61").Value = 67430.664
$ws.Range("L61").Value = 113110.336
$ws.Range("M61").Value = -67228.664
$ws.Range("N61").Value = -113514.336
$ws.Range("H99").Value = 18820.75
$ws.Range("I99").Value = 14999.333
$ws.Range("J99").Value = 30285
$ws.Range("K99").Value = 14999.333
$ws.Range("L99").Value = 30285
$ws.Range("M99").Value = -12004.333
$ws.Range("N99").Value = -36275
$ws.Range("H113").Value = 75043.94500000001
$ws.Range("I113").Value = 67430.664
$ws.Range("J113").Value = 113110.336
$ws.Range("K113").Value = 67430.664
$ws.Range("L113").Value = 113110.336
$ws.Range("M113").Value = -65260.664
$ws.Range("N113").Value = -117450.336
$ws.Range("H126").Value = 19791.438
$ws.Range("I126").Value = 35717.875
$ws.Range("K126").Value = 107153.625
$ws.Range("M126").Value = -104683.625
$ws.Range("H132").Value = 2572.7192
$ws.Range("I132").Value = 2040.8605
$ws.Range("K132").Value = 6122.5815
$ws.Range("M132").Value = -3592.5815
$ws.Range("H136").Value = 2303.4468
$ws.Range("I136").Value = 1938.8206
$ws.Range("J136").Value = 4081
$ws.Range("K136").Value = 5816.4618
$ws.Range("L136").Value = 12243
$ws.Range("M136").Value = -3266.4618
$ws.Range("N136").Value = -17343

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4287.2383
$ws.Range("I81").Value = 4235.1113
$ws.Range("K81").Value = 8470.222599999999
$ws.Range("M81").Value = -7409.222599999999
$ws.Range("H84").Value = 4287.2383
$ws.Range("I84").Value = 4235.1113
$ws.Range("K84").Value = 42351.113
$ws.Range("M84").Value = -37047.113
$ws.Range("H100").Value = 2095.2778
$ws.Range("I100").Value = 2205.4707
$ws.Range("K100").Value = 4410.9414
$ws.Range("M100").Value = -3869.9414
